# Update the "bugs" sheet of the bug-report workbook:
#  - bug #8 (row 11) status changes from "open" to "won't fix (this is Windows behaviour)"
#  - bug #9 (row 12) status changes from "open" to "resolved"
#  - widen the "Status" column (E) a bit
#  - leave the view scrolled/selected near the bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bugs")

# -- content updates (Status column) --------------------------------------
$ws.Range("E11").Value = "won't fix (this is Windows behaviour)"
$ws.Range("E12").Value = "resolved"

# -- formatting: widen column E (Status) -----------------------------------
$ws.Columns.Item(5).ColumnWidth = 34.6666666666667

# -- view state: scroll down a bit and move the selection -------------------
$ws.Activate()
$ws.Range("G14").Select()
